$d = $word.ActiveDocument

# Replace the old answer "|  (pipe symbol)" with ">" for the first question.
$d.Content.Find.Execute("|  (pipe symbol)", $false, $false, $false, $false, $false, $true, 1, $false, ">", 2)
